# Re-run SGNN to annotate dialog acts following clean up work to the original transcripts.
# Update DAMSLTag (column I) and DialogAct (column J) values for the affected rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 4;   Tag = "sv"; Act = "Statement-opinion" },
    @{ Row = 10;  Tag = "sv"; Act = "Statement-opinion" },
    @{ Row = 31;  Tag = "aa"; Act = "Agree/Accept" },
    @{ Row = 74;  Tag = "sd"; Act = "Statement-non-opinion" },
    @{ Row = 75;  Tag = "ba"; Act = "Appreciation" },
    @{ Row = 78;  Tag = "ba"; Act = "Appreciation" },
    @{ Row = 85;  Tag = "sd"; Act = "Statement-non-opinion" },
    @{ Row = 92;  Tag = "sd"; Act = "Statement-non-opinion" },
    @{ Row = 116; Tag = "sv"; Act = "Statement-opinion" },
    @{ Row = 121; Tag = "b";  Act = "Acknowledge (Backchannel)" },
    @{ Row = 123; Tag = "%";  Act = "Uninterpretable" }
)

foreach ($u in $updates) {
    $ws.Range("I$($u.Row)").Value = $u.Tag
    $ws.Range("J$($u.Row)").Value = $u.Act
}
